{"js": "// Replace the bold \"Concept 3:\" paragraph with the new concluding sentence,\n// clear out the unfinished \"Step 1. Man to \" paragraph, and drop the\n// \"_GoBack\" bookmark at the very end of the document (marking the last\n// edit position), matching Word's behavior when a document is saved.\n\nconst body = context.document.body;\n\n// 1) \"Concept 3:\" -> \"Both concepts are viable, however there is less risk\n//    for the animals in Concept 2. \" (not bold, unlike the old heading text).\nconst conceptResults = body.search(\"Concept 3:\", { matchCase: true });\nconceptResults.load(\"items\");\nawait context.sync();\n\nconst conceptRange = conceptResults.items[0];\nconceptRange.font.bold = false;\nconceptRange.insertText(\n  \"Both concepts are viable, however there is less risk for the animals in Concept 2. \",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// 2) Remove the unfinished \"Step 1. Man to \" text, leaving the paragraph empty.\nconst stepResults = body.search(\"Step 1. Man to \", { matchCase: true });\nstepResults.load(\"items\");\nawait context.sync();\n\nconst stepRange = stepResults.items[0];\nstepRange.insertText(\"\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 3) Leave a \"_GoBack\" bookmark at the last paragraph (the final empty\n//    paragraph at the end of the document), as Word does for the most\n//    recent edit location.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.getRange().insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Replace the bold \"Concept 3:\" paragraph with the new concluding sentence,\n# clear out the unfinished \"Step 1. Man to \" paragraph, and drop a\n# \"_GoBack\" bookmark at the very end of the document (marking the last\n# edit position), matching Word's behavior when a document is saved.\n\n$d = $word.ActiveDocument\n\n# 1) \"Concept 3:\" -> \"Both concepts are viable, however there is less risk\n#    for the animals in Concept 2. \" (not bold, unlike the old heading text).\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Execute(\"Concept 3:\") | Out-Null\n$rng.Text = \"Both concepts are viable, however there is less risk for the animals in Concept 2. \"\n$rng.Font.Bold = 0\n\n# 2) Remove the unfinished \"Step 1. Man to \" text, leaving the paragraph empty.\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Execute(\"Step 1. Man to \") | Out-Null\n$rng2.Text = \"\"\n\n# 3) Leave a \"_GoBack\" bookmark at the last paragraph (the final empty\n#    paragraph at the end of the document), as Word does for the most\n#    recent edit location.\n$lastParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$d.Bookmarks.Add(\"_GoBack\", $lastParagraph.Range) | Out-Null\n"}
